$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) format, used to restore style
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '43.429.97'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").Value = '2.376.77'
$ws.Range("E3").Value = '  +5.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $normalStyle
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.87'
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.657'
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = '  +4.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.94'
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = '  +14.59%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.465'
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = '  +3.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0974'
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.07'
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.52'
$ws.Range("D12").Style = $normalStyle
$ws.Range("D13").Value = '2.727.08'
$ws.Range("E13").Value = '  +5.61%  '
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.81'
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.26'
$ws.Range("D16").Style = $normalStyle
$ws.Range("E16").Value = '  +3.53%  '
$ws.Range("E17").Value = '  +3.53%  '
$ws.Range("D18").Value = '2.375.52'
$ws.Range("E18").Value = '  +5.80%  '
$ws.Range("D19").Value = '43.438.91'
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("D20").Value = '0.0₃0991'
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.38'
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = '  +5.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.40'
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = '  +2.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.56'
$ws.Range("D23").Style = $normalStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.98'
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = '  +19.88%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.48'
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = '  +2.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.47'
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = '  +11.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = '  +2.65%  '
$ws.Range("E29").Value = '  -6.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.42'
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("E31").Value = '  +9.01%  '
$ws.Range("E32").Value = '  -8.75%  '
$ws.Range("E33").Value = '  +2.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.01'
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = '  +4.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0694'
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = '  +1.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.11'
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.46'
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = '  +8.43%  '
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.60'
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = '  +3.40%  '
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0255'
$ws.Range("D40").Style = $normalStyle
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.96'
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = '  +3.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.67'
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = '  +9.59%  '
$ws.Range("E44").Value = '  +10.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.73'
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = '  +2.28%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.22'
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = '  +2.40%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.50'
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = '  +3.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0950'
$ws.Range("D48").Style = $normalStyle
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("D49").Value = '1.454.47'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D50").Value = '2.599.38'
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.74'
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = '  -0.60%  '
